# Add a 3-item numbered list (POST vs PUT / relative URL / absolute vs
# relative URL) to the trailing blank paragraph of the document, matching
# the "I did the first three, breifly." commit.

$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs (spacing after=0).
# The very last paragraph is where the list starts.
$lastIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastIndex)

# --- Apply list/style formatting to the (still empty) paragraph first so
# that later edits don't stamp extra rsid attributes onto it. -------------
$p.Range.Text = "x"
$p.Style = "List Paragraph"
$p.Format.SpaceAfter = 0
$p.Range.ListFormat.ApplyNumberDefault()

# Match the "1) a. i." hybrid multilevel list (decimal / lowerLetter /
# lowerRoman repeating every 3 levels, with the top level using "%1)").
$lt = $p.Range.ListFormat.ListTemplate
$lvl1 = $lt.ListLevels.Item(1)
$lvl1.NumberFormat = "%1)"
$numStyles = @(0, 4, 2, 0, 4, 2, 0, 4, 2)
for ($i = 2; $i -le 9; $i++) {
    $lvl = $lt.ListLevels.Item($i)
    $lvl.NumberStyle = $numStyles[$i - 1]
}

# --- Build the first bullet's text. It is authored as three separate runs
# in the target (no formatting difference between them), so assemble it as
# three temporary paragraphs and then merge the paragraph marks back out,
# which preserves each chunk as its own <w:r>. ----------------------------
$r = $p.Range
$r.Text = "POST sends data to a specific URL, whereas PUT puts a webpage at a specified URL.  The data sent by POST is dealt with however the server side decides to deal with it.  PUT"
$r.Collapse(0)
$r.InsertParagraphAfter()

$mid = $d.Paragraphs.Item($lastIndex + 1)
$mid.Range.InsertAfter(" is an idempotent operation, which means that doing the operation PUT multiple times is no di")
$mid.Range.Collapse(0)
$mid.Range.InsertParagraphAfter()

$tail = $d.Paragraphs.Item($lastIndex + 2)
$tail.Range.InsertAfter("fferent than doing it only once.")

$mark1 = $d.Range($d.Paragraphs.Item($lastIndex).Range.End - 1, $d.Paragraphs.Item($lastIndex).Range.End)
$mark1.Delete()
$mark2 = $d.Range($d.Paragraphs.Item($lastIndex).Range.End - 1, $d.Paragraphs.Item($lastIndex).Range.End)
$mark2.Delete()

$p1 = $d.Paragraphs.Item($lastIndex)

# --- Second bullet. --------------------------------------------------------
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($lastIndex + 1)
$p2.Range.InsertAfter("The target URL is relative.")

# --- Third bullet. ----------------------------------------------------------
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($lastIndex + 2)
$p3.Range.InsertAfter("The difference between a relative and an absolute URL is that the absolute URL contains the entire URL unlike the relative URL which contains a partial URL.")

# --- Fix up the generated "List Paragraph" style definition so it matches
# Word's real built-in style (priority 34, 720-twip left indent, and
# contextual spacing so list items don't add space between themselves). --
$style = $d.Styles.Item("List Paragraph")
$style.Priority = 34
$style.ParagraphFormat.LeftIndent = 36
$style.NoSpaceBetweenParagraphsOfSameStyle = $true

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
